$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-run-of-model numeric updates in the capital cost block (rows 3-17)
# ---------------------------------------------------------------------------
$ws.Range("C3").Value  = 65.787766180879586
$ws.Range("C4").Value  = 32.116801970507218
$ws.Range("C5").Value  = 2.6315106472351841
$ws.Range("C6").Value  = 5.9208989562791627
$ws.Range("C7").Value  = 2.9604494781395809
$ws.Range("C8").Value  = 109.4174272330407
$ws.Range("C9").Value  = 10.94174272330407
$ws.Range("C10").Value = 10.94174272330407
$ws.Range("C11").Value = 21.883485446608152
$ws.Range("C12").Value = 32.825228169912222
$ws.Range("C13").Value = 10.94174272330407
$ws.Range("C14").Value = 87.533941786432621
$ws.Range("C15").Value = 196.95136901947339
$ws.Range("C16").Value = 9.8475684509736681
$ws.Range("C17").Value = 206.79893747044699

# ---------------------------------------------------------------------------
# 2) Variable operating-cost block (rows 21-28) is reordered (raw materials
#    regrouped: DAP, Salt / Wastewater / Glucose, Process water, Tridecane,
#    CSL, Natural gas) and the merged "Raw materials" label now spans two
#    separate groups (A21:A22 and A24:A28), with "By-products and credits"
#    moving up to row 23 on its own.
# ---------------------------------------------------------------------------

# Break the old merges first so new ones can be created cleanly.
$ws.Range("A21:A25").UnMerge()
$ws.Range("A27:A28").UnMerge()

# Row 21 - Raw materials / DAP
$ws.Range("A21").Value2 = "Raw materials"
$ws.Range("B21").Value2 = "DAP"
$ws.Range("C21").Value  = 895.39159499999994
$ws.Range("D21").Value  = 1.1811074895224249

# Row 22 - (merged with 21) / Salt
$ws.Range("A22").Value2 = ""
$ws.Range("B22").Value2 = "Salt"
$ws.Range("C22").Value  = 136.07775000000001
$ws.Range("D22").Value  = 4.2304986641846387

# Row 23 - By-products and credits / Wastewater
$ws.Range("A23").Value2 = "By-products and credits"
$ws.Range("B23").Value2 = "Wastewater"
$ws.Range("C23").Value  = -1.9630611699040379
$ws.Range("D23").Value  = -6.2658539245062812

# Row 24 - Raw materials / Glucose
$ws.Range("A24").Value2 = "Raw materials"
$ws.Range("B24").Value2 = "Glucose"
$ws.Range("C24").Value  = 244.03276500000001
$ws.Range("D24").Value  = 146.92179801857401

# Row 25 - (merged with 24) / Process water
$ws.Range("A25").Value2 = ""
$ws.Range("B25").Value2 = "Process water"
$ws.Range("C25").Value  = 0.320236305
$ws.Range("D25").Value  = 0.99534905024449527

# Row 26 - (merged with 24) / Tridecane
$ws.Range("A26").Value2 = ""
$ws.Range("B26").Value2 = "Tridecane"
$ws.Range("C26").Value  = 878.15507999999988
$ws.Range("D26").Value  = 0.0024024701315892168

# Row 27 - (merged with 24) / CSL
$ws.Range("A27").Value2 = ""
$ws.Range("B27").Value2 = "CSL"
$ws.Range("C27").Value  = 51.528108000000003
$ws.Range("D27").Value  = 0.56651201273216611

# Row 28 - (merged with 24) / Natural gas
$ws.Range("A28").Value2 = ""
$ws.Range("B28").Value2 = "Natural gas"
$ws.Range("C28").Value  = 197.76633000000001
$ws.Range("D28").Value  = 3.8546465757895638

# Re-create the merges in their new groupings.
$ws.Range("A21:A22").Merge()
$ws.Range("A24:A28").Merge()

# Give the (new) "Raw materials" label cells their own clone of the bordered
# header style (matches workbook's duplicated cellXfs entry for this block).
$ws.Range("A21:A22").Font.Bold = $true
$ws.Range("A21:A22").Borders.LineStyle = 1
$ws.Range("A21:A22").HorizontalAlignment = -4108
$ws.Range("A21:A22").VerticalAlignment = -4160

$ws.Range("A24:A28").Font.Bold = $true
$ws.Range("A24:A28").Borders.LineStyle = 1
$ws.Range("A24:A28").HorizontalAlignment = -4108
$ws.Range("A24:A28").VerticalAlignment = -4160

# Total variable operating cost
$ws.Range("D29").Value = 164.01816820568561

# ---------------------------------------------------------------------------
# 3) Labor / other cost block (rows 35-36) re-run values
# ---------------------------------------------------------------------------
$ws.Range("C35").Value = 1.973632985426387
$ws.Range("D35").Value = 1.8946876660093319
$ws.Range("C36").Value = 0.46051436326615708
$ws.Range("D36").Value = 0.44209378873551092

# ---------------------------------------------------------------------------
# 4) Column widths / view (cosmetic, matches a resave in a newer Excel)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37.88671875
$ws.Columns.Item(2).ColumnWidth = 16.6640625
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 14.5546875

$ws.Range("D24").Select()
